$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("D5").Value = "récupère la valeur du champ ""quantité"" sur la page product afin d'etre après ajoutée au localStorage"
$ws.Range("E5").Value = "au bouton ajouter au panier, il se retrouve dans le localStorage"
$ws.Range("F5").Value = "ajout a zero ou si la couleur n'a pas été choisie doit etre évitée"

# Row 6
$ws.Range("A6").Value = "product.js"
$ws.Range("B6").Value = "50 à 53"
$ws.Range("C6").Value = "colorValue()"
$ws.Range("D6").Value = "récupère la valeur du champ ""couleur"" sur la page product afin d'etre après ajoutée au localStorage"
$ws.Range("E6").Value = "au bouton ajouter au panier, il se retrouve dans le localStorage"
$ws.Range("F6").Value = "si ""aucune couleur selectionnée"" devrait ne pas fonctionner, mais peut etre un probleme"

# Row 7
$ws.Range("A7").Value = "product.js"
$ws.Range("B7").Value = "58 à 62"
$ws.Range("C7").Value = "toCartBtn.addEventListener(""click"", () "
$ws.Range("D7").Value = "récupère les valeurs de ""qtyValue() et colorValue()"" et les envoie a la fction add2cart() de la page cart.js avec l'id du produit définie dans les 3 premieres lignes de product.js"
$ws.Range("E7").Value = "au bouton ajouter au panier, id, color et qty se retrouvent dans le localStorage"
$ws.Range("F7").Value = "si une des trois valeur, ""id, color ou qty"" est nulle cela devrait etre incomplet dans le localStorage"

# Row 8
$ws.Range("A8").Value = "cart.js"
$ws.Range("B8").Value = "6 à 12"
$ws.Range("C8").Value = "getCart()"
$ws.Range("D8").Value = "récupère la valeur ""panier"" du localStorage dans une variable items qui est la sortie de la fonction"
$ws.Range("E8").Value = "utilisé dans plusieurs fonctions, la fonction getCart permet de stocker dans une variable le localStorage, notamment dans ""add2cart(), deleteItem(), changeQuantity(), fetchIdData(), makeJsonData()"". Cette variable est tout le temps items donc il suffit de surveiller items dans l'inspecteur"
$ws.Range("F8").Value = "si le panier dans le localStorage n'est pas complet, items ne sera pas complet et tout bug. Mais cette fonction tres utilisée est essentielle a d'autres fonctions donc les bugs peuvent etre multiple si elle ne fonctionne pas"

# Row 9
$ws.Range("A9").Value = "cart.js"
$ws.Range("B9").Value = "15 à 36"
$ws.Range("C9").Value = "add2cart()"
$ws.Range("D9").Value = "construit le panier du localStorage, prend en compte le panier existant"
$ws.Range("E9").Value = "regarder le localStorage"
$ws.Range("F9").Value = "definitions de items par getCart() mauvaise ; ou envoie de ""id, color et qty""par ""toCartBtn.addEventListener(""click"", () "" dans la page product mauvais aussi."

# Row 10
$ws.Range("A10").Value = "cart.js"
$ws.Range("B10").Value = "39 à 52"
$ws.Range("C10").Value = "deleteItem()"
$ws.Range("D10").Value = "supprime une entrée du panier de html et du localStorage, et reload la page"
$ws.Range("E10").Value = "l'entrée html doit etre supprimée et l'entrée localStorage doit être supprimée"
$ws.Range("F10").Value = "la page pourrait ne pas se recharger, et l'html ne pas effacer la zone du kanap a supprimer. ""splice()"" pourrait rendre un mauvais ""items"" et dans ce cas corrompre le nouveau panier"

# Row 11
$ws.Range("A11").Value = "cart.js"
$ws.Range("B11").Value = "54 à 63"
$ws.Range("C11").Value = "changeQuantity()"
$ws.Range("D11").Value = "modifie la quantité d'un item demandé par l'utilisateur dans le local storage."
$ws.Range("E11").Value = "Lorsque on modifie la quantité d'un kanap dans la panier, le localstorage doit etre immediatement modifié pour la valeur du client."
$ws.Range("F11").Value = "Si ""getCart()"" est corrompu la suite est corrompue. Si le panier ne correspond pas au localstorage, le changement ne se fera pas"

# Row 12
$ws.Range("A12").Value = "cart.js"
$ws.Range("B12").Value = "71 à 118"
$ws.Range("C12").Value = "fetchIdData()"
$ws.Range("D12").Value = "affiche en innerHTML les éléments du panier et leur contenu de l'API. Affiche panier vide si panier vide."
$ws.Range("E12").Value = "les éléménets s'affichent correctement et completement : l'image, son texte alt, le nom du kanap, sa couleur, son prix, sa quantité désirée, son prix total, le nombre d'articles total"
$ws.Range("F12").Value = "Vérifier que le prix total et le nombre total d'article soit le bon. Problemes de fetch de l'API, si le back end n'est pas allumé par exemple. ""getCart()"" peut etre corrompu. La panier pourrait ne pas s'afficher comme vide en HTML dans le cas ou localStorage est null"

# Row 13
$ws.Range("A13").Value = "cart.js"
$ws.Range("B13").Value = "133 à 142"
$ws.Range("C13").Value = "validateEmail(mai)"
$ws.Range("D13").Value = "Fonction REGEX pour valider une adresse email"
$ws.Range("E13").Value = "si l'utilisateur entre une adresse non conforme, cette fonction affiche false. True pour l'inverse. Elle est declenchée par ""orderButton.addEventListener(""click"", (e) =>"""
$ws.Range("F13").Value = "la regex pourrait laisser passer de mauvaises adresses mail."

# Row 14
$ws.Range("A14").Value = "cart.js"
$ws.Range("B14").Value = "147 à 154"
$ws.Range("C14").Value = "validateFirstName(prenom)"
$ws.Range("D14").Value = "Fonction regex pour valider un prenom sans chiffre"
$ws.Range("E14").Value = "cette fonction ne doit laisser passer aucun chiffre"
$ws.Range("F14").Value = "la regex pourrait laisser passer des prenoms invalides"

# Row 15
$ws.Range("A15").Value = "cart.js"
$ws.Range("B15").Value = "158 à 165"
$ws.Range("C15").Value = "validateLastName(nom)"
$ws.Range("D15").Value = "Fonction regex pour valider un nom sans chiffre"
$ws.Range("E15").Value = "cette fonction ne doit laisser passer aucun chiffre"
$ws.Range("F15").Value = "la regex pourrait laisser passer des noms invalides"

# Row 16
$ws.Range("A16").Value = "cart.js"
$ws.Range("B16").Value = "169 à 176"
$ws.Range("C16").Value = "validateCity(ville)"
$ws.Range("D16").Value = "Fonction regex pour valider une ville sans chiffre"
$ws.Range("E16").Value = "cette fonction ne doit laisser passer aucun chiffre"
$ws.Range("F16").Value = "la regex pourrait laisser passer des villes invalides"

# Row 17
$ws.Range("A17").Value = "cart.js"
$ws.Range("B17").Value = "181 à 188"
$ws.Range("C17").Value = "validateAddress(adresse)"
$ws.Range("D17").Value = "Fonction regex pour valider une adresse non vide"
$ws.Range("E17").Value = "cette fonction demande au champ d'etre non-vide"
$ws.Range("F17").Value = "la regex pourrait laisser passer des champs vides"

# Row 18
$ws.Range("A18").Value = "cart.js"
$ws.Range("B18").Value = "208 à 231"
$ws.Range("C18").Value = "makeJsonData()"

# Row 19
$ws.Range("A19").Value = "cart.js"
$ws.Range("B19").Value = "235 à 285"
$ws.Range("C19").Value = "orderButton.addEventListener(""click"", (e) =>"

# Update selection to reflect final active cell per the authored edit
$ws.Range("B19").Select()
